# Refresh the scraped crypto table (Coin/Link/Price/Volume) to match the
# Jan 6 2024 GitHub Actions data pull: per-row Price (D) / Volume 1h (E) updates,
# plus a handful of rows whose rank reshuffled so Coin (B) and Link (C) changed too.
#
# All of these cells are plain text in the workbook (e.g. D holds "1.00",
# "0.0810", "44.245.45" as literal strings, not numbers). Range.Value2 keeps a
# string as a string UNLESS it parses cleanly as a number, in which case COM
# auto-converts it like typing it into a cell would - so for those specific
# values we force the Text number format first to keep them as exact strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value2 = '44.245.45'
$ws.Range("E2").Value2 = '  +0.59%  '

# Row 3: update D3, E3
$ws.Range("D3").Value2 = '2.239.89'
$ws.Range("E3").Value2 = '  -0.03%  '

# Row 4: update D4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.84'
$ws.Range("E5").Value2 = '  -2.48%  '

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.04'
$ws.Range("E6").Value2 = '  -5.08%  '

# Row 7: update E7
$ws.Range("E7").Value2 = '  -0.51%  '

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value2 = '  +0.25%  '

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value2 = '  -1.22%  '

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.62'
$ws.Range("E10").Value2 = '  -4.00%  '

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value2 = '  -1.36%  '

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.17'
$ws.Range("E12").Value2 = '  -2.29%  '

# Row 13: update E13
$ws.Range("E13").Value2 = '  -0.17%  '

# Row 14: update B14, C14, D14, E14
$ws.Range("B14").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value2 = '2.581.24'
$ws.Range("E14").Value2 = '  +0.03%  '

# Row 15: update B15, C15, D15, E15
$ws.Range("B15").Value2 = 'WrappedEther'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value2 = '2.240.73'
$ws.Range("E15").Value2 = '  +0.22%  '

# Row 16: update D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.829'
$ws.Range("E16").Value2 = '  -1.22%  '

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.49'
$ws.Range("E17").Value2 = '  -2.98%  '

# Row 18: update D18, E18
$ws.Range("D18").Value2 = '43.959.75'
$ws.Range("E18").Value2 = '  +0.31%  '

# Row 19: update D19, E19
$ws.Range("D19").Value2 = '0.0₃0963'
$ws.Range("E19").Value2 = '  -1.31%  '

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.38'
$ws.Range("E20").Value2 = '  +1.09%  '

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.09'
$ws.Range("E21").Value2 = '  -7.56%  '

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.55'
$ws.Range("E22").Value2 = '  +0.14%  '

# Row 23: update B23, C23, D23, E23
$ws.Range("B23").Value2 = 'BitcoinCash'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.62'
$ws.Range("E23").Value2 = '  +0.85%  '

# Row 24: update B24, C24, D24, E24
$ws.Range("B24").Value2 = 'PancakeSwap'
$ws.Range("C24").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.94'
$ws.Range("E24").Value2 = '  -1.16%  '

# Row 25: update E25
$ws.Range("E25").Value2 = '  -0.93%  '

# Row 26: update E26
$ws.Range("E26").Value2 = '  +0.22%  '

# Row 27: update B27, C27, D27, E27
$ws.Range("B27").Value2 = 'Toncoin'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.22'
$ws.Range("E27").Value2 = '  +4.12%  '

# Row 28: update B28, C28, D28, E28
$ws.Range("B28").Value2 = 'InjectiveProtocol'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.27'
$ws.Range("E28").Value2 = '  +5.51%  '

# Row 29: update B29, C29, D29, E29
$ws.Range("B29").Value2 = 'Cosmos'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.89'
$ws.Range("E29").Value2 = '  -2.45%  '

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.02'
$ws.Range("E30").Value2 = '  +0.00%  '

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.84'
$ws.Range("E31").Value2 = '  -2.01%  '

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.77'
$ws.Range("E32").Value2 = '  -1.54%  '

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0794'
$ws.Range("E33").Value2 = '  -4.83%  '

# Row 34: update E34
$ws.Range("E34").Value2 = '  -1.82%  '

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.14'
$ws.Range("E35").Value2 = '  -4.36%  '

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("E36").Value2 = '  +1.97%  '

# Row 37: update E37
$ws.Range("E37").Value2 = '  -1.01%  '

# Row 38: update E38
$ws.Range("E38").Value2 = '  -7.23%  '

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.50'
$ws.Range("E39").Value2 = '  +0.00%  '

# Row 40: update E40
$ws.Range("E40").Value2 = '  -4.67%  '

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.32'
$ws.Range("E41").Value2 = '  -7.33%  '

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0299'
$ws.Range("E42").Value2 = '  -2.28%  '

# Row 43: update E43
$ws.Range("E43").Value2 = '  +0.24%  '

# Row 44: update D44, E44
$ws.Range("D44").Value2 = '1.744.63'
$ws.Range("E44").Value2 = '  +2.81%  '

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.70'
$ws.Range("E45").Value2 = '  +0.33%  '

# Row 46: update E46
$ws.Range("E46").Value2 = '  -1.96%  '

# Row 47: update B47, C47, D47, E47
$ws.Range("B47").Value2 = 'Aave'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.69'
$ws.Range("E47").Value2 = '  -1.81%  '

# Row 48: update B48, C48, D48, E48
$ws.Range("B48").Value2 = 'THORChain'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.92'
$ws.Range("E48").Value2 = '  -4.74%  '

# Row 49: update B49, C49, D49, E49
$ws.Range("B49").Value2 = 'FraxShare'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.06'
$ws.Range("E49").Value2 = '  -0.67%  '

# Row 50: update B50, C50, D50, E50
$ws.Range("B50").Value2 = 'Stacks'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.57'
$ws.Range("E50").Value2 = '  -1.37%  '

# Row 51: update B51, C51, D51, E51
$ws.Range("B51").Value2 = 'MultiversX'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.44'
$ws.Range("E51").Value2 = '  -2.91%  '
